$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 14-20) so the sheet's
# dimension shrinks back down to A1:F13.
$ws.Range("A14:A20").EntireRow.Delete()

# New data for rows 2-13 (years 2010-2021), shifted up from the old rows
# 10-20 plus one new row (2021) appended at the bottom.
$data = @(
    @("2010年", 218.0577, 10.4,               56.07,              3.47404131, 1.37111392),
    @("2011年", 247.3426, 10.1506614516,      54.3534497504,      4.09499505, 1.37037903),
    @("2012年", 268.6554, 10.0544507384,      55.4717994843,      4.54751077, 1.43935953),
    @("2013年", 292.063,  9.8348006948,       57.0199164976,      5.07885866, 1.49211933),
    @("2014年", 298.0571, 9.858046419900001,  55.6220400868,      5.36187933, 1.49119766),
    @("2015年", 305.5499, 9.800000000000001,  54.7,               5.59026,    1.47425),
    @("2016年", 313.7143, 9.704093051799999,  54.6218001719,      5.63270221, 1.55618949),
    @("2017年", 344.2497, 9.5,                54.84,              6.07432288, 1.59823646),
    @("2018年", 339.5371, 9.9,                52,                 6.38978662, 1.60115334),
    @("2019年", 339.52,   9.699999999999999,  49.65,              6.911069,   1.680566),
    @("2020年", 292.7288, 10.3,               42.81,              6.20683853, 1.34037234),
    @("2021年", 319.2971, 9.800000000000001,  43.2,               6.95966057, 1.40058791)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row = $row + 1
}
